# Generate Report for Handoff
# Adds two new localization entries:
#   - 3414b3fb-bece-4317-851b-2080f8fb27f6.md  (new row, inserted before 7f880207 entry)
#   - f5fea1db-e187-4e4e-a926-19f06f96b483.md  (new row, appended at the end)
# across the Overview, zh-cn and de-de sheets/tables.

$wb = $excel.ActiveWorkbook

$hlColor = 15570276   # RGB(0x64,0x95,0xED) == ARGB FF6495ED used by the workbook's hyperlink style

function Style-Hyperlink($range) {
    $range.Font.Underline = $true
    $range.Font.Color = $hlColor
}

function Style-DateText($range) {
    $range.NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

# ---------------------------------------------------------------------------
# Sheet "Overview" (columns: File Name | Path And Name | Extension | Publish URL | zh-cn | de-de | Latest HO Xliff Generate Date)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Row 5 becomes the new 3414b3fb entry
$wsOverview.Range("A5").Value = "3414b3fb-bece-4317-851b-2080f8fb27f6.md"
$wsOverview.Range("B5").Value = "e2e\3414b3fb-bece-4317-851b-2080f8fb27f6.md"
$wsOverview.Range("C5").Value = ".md"
$wsOverview.Range("E5").Value = "Ready for handoff"
$wsOverview.Range("F5").Value = "Ready for handoff"
$wsOverview.Range("G5").Value = "2016-08-31 12:51:36"
Style-DateText $wsOverview.Range("G5")

# Row 6: the previous 7f880207 entry, shifted down
$wsOverview.Range("A6").Value = "7f880207-dc89-446c-99de-d3b10383c0e3.md"
$wsOverview.Range("B6").Value = "e2e\7f880207-dc89-446c-99de-d3b10383c0e3.md"
$wsOverview.Range("C6").Value = ".md"
$wsOverview.Range("E6").Value = "Ready for handoff"
$wsOverview.Range("F6").Value = "Ready for handoff"
$wsOverview.Range("G6").Value = "2016-08-31 12:49:20"
Style-DateText $wsOverview.Range("G6")

# Row 7: new f5fea1db entry
$wsOverview.Range("A7").Value = "f5fea1db-e187-4e4e-a926-19f06f96b483.md"
$wsOverview.Range("B7").Value = "e2e\f5fea1db-e187-4e4e-a926-19f06f96b483.md"
$wsOverview.Range("C7").Value = ".md"
$wsOverview.Range("E7").Value = "Ready for handoff"
$wsOverview.Range("F7").Value = "Ready for handoff"
$wsOverview.Range("G7").Value = "2016-08-31 12:51:36"
Style-DateText $wsOverview.Range("G7")

# Hyperlinks on column B (display text already set above, so pass the existing text through)
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/90e4f0ebcb4243bb0e92ca46307a9286e675fe81/e2e/e75892d4-044b-4d09-8abf-855ed5fa4f20.md", "", "", "e2e\e75892d4-044b-4d09-8abf-855ed5fa4f20.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cd76318f29306b8e90dd81f65ec88bd3a309cc5e/e2e/7c27a539-13f8-415c-9532-2e5a2ec12fb6.md", "", "", "e2e\7c27a539-13f8-415c-9532-2e5a2ec12fb6.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cd76318f29306b8e90dd81f65ec88bd3a309cc5e/e2e/acfcdcf1-9169-4293-80a2-ed8587dc1452.md", "", "", "e2e\acfcdcf1-9169-4293-80a2-ed8587dc1452.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d92d0733590fd5894259fe924c051e691cdc04b6/e2e/3414b3fb-bece-4317-851b-2080f8fb27f6.md", "", "", "e2e\3414b3fb-bece-4317-851b-2080f8fb27f6.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d92d0733590fd5894259fe924c051e691cdc04b6/e2e/7f880207-dc89-446c-99de-d3b10383c0e3.md", "", "", "e2e\7f880207-dc89-446c-99de-d3b10383c0e3.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d92d0733590fd5894259fe924c051e691cdc04b6/e2e/f5fea1db-e187-4e4e-a926-19f06f96b483.md", "", "", "e2e\f5fea1db-e187-4e4e-a926-19f06f96b483.md") | Out-Null

Style-Hyperlink $wsOverview.Range("B2")
Style-Hyperlink $wsOverview.Range("B3")
Style-Hyperlink $wsOverview.Range("B4")
Style-Hyperlink $wsOverview.Range("B5")
Style-Hyperlink $wsOverview.Range("B6")
Style-Hyperlink $wsOverview.Range("B7")

# Resize the "Overview" table to include the two new rows
$tblOverview = $wsOverview.ListObjects.Item(1)
$tblOverview.Resize($wsOverview.Range("A1:G7"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# (columns: A Source File Name, B File Extension, C Status, D Source Path,
#  E Priority, F Content Duplicate, G Latest Handoff File, H Latest Handoff Datetime,
#  I Latest Target File, J Latest Handback File, K Latest Handback DateTime,
#  L Reference Tokens, M To be localized, N Dependency From, O Has metadata, P Error Detail)
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# Row 5: new 3414b3fb entry
$wsZhCn.Range("A5").Value = "3414b3fb-bece-4317-851b-2080f8fb27f6.md"
$wsZhCn.Range("B5").Value = ".md"
$wsZhCn.Range("C5").Value = "Ready for handoff"
$wsZhCn.Range("D5").Value = "e2e"
$wsZhCn.Range("E5").Value = "ht"
$wsZhCn.Range("F5").Value = "False"
$wsZhCn.Range("G5").Value = "3414b3fb-bece-4317-851b-2080f8fb27f6.9a7d471b31670a0c41801e2ea85158d6511dc541.zh-cn.xlf"
$wsZhCn.Range("H5").Value = "2016-08-31 12:51:32"
$wsZhCn.Range("K5").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("M5").Value = "True"
$wsZhCn.Range("O5").Value = "False"
Style-DateText $wsZhCn.Range("H5")
Style-DateText $wsZhCn.Range("K5")

# Row 6: previous 7f880207 entry, shifted down
$wsZhCn.Range("A6").Value = "7f880207-dc89-446c-99de-d3b10383c0e3.md"
$wsZhCn.Range("B6").Value = ".md"
$wsZhCn.Range("C6").Value = "Ready for handoff"
$wsZhCn.Range("D6").Value = "e2e"
$wsZhCn.Range("E6").Value = "ht"
$wsZhCn.Range("F6").Value = "False"
$wsZhCn.Range("G6").Value = "7f880207-dc89-446c-99de-d3b10383c0e3.fc6f83bef6d89ff41f4a184c043c9e87ea8c1cdf.zh-cn.xlf"
$wsZhCn.Range("H6").Value = "2016-08-31 12:49:15"
$wsZhCn.Range("K6").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("M6").Value = "True"
$wsZhCn.Range("O6").Value = "False"
Style-DateText $wsZhCn.Range("H6")
Style-DateText $wsZhCn.Range("K6")

# Row 7: new f5fea1db entry
$wsZhCn.Range("A7").Value = "f5fea1db-e187-4e4e-a926-19f06f96b483.md"
$wsZhCn.Range("B7").Value = ".md"
$wsZhCn.Range("C7").Value = "Ready for handoff"
$wsZhCn.Range("D7").Value = "e2e"
$wsZhCn.Range("E7").Value = "ht"
$wsZhCn.Range("F7").Value = "False"
$wsZhCn.Range("G7").Value = "f5fea1db-e187-4e4e-a926-19f06f96b483.9e74bf257fb01b4211733f5bb469e0eb170f6bfd.zh-cn.xlf"
$wsZhCn.Range("H7").Value = "2016-08-31 12:51:32"
$wsZhCn.Range("K7").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("M7").Value = "True"
$wsZhCn.Range("O7").Value = "False"
Style-DateText $wsZhCn.Range("H7")
Style-DateText $wsZhCn.Range("K7")

# Hyperlinks: column A for every data row, plus the extra one on I2
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/90e4f0ebcb4243bb0e92ca46307a9286e675fe81/e2e/e75892d4-044b-4d09-8abf-855ed5fa4f20.md", "", "", "e75892d4-044b-4d09-8abf-855ed5fa4f20.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/c7d36a619a7c8394eebf65c63146980dabae39c9/e2e/e75892d4-044b-4d09-8abf-855ed5fa4f20.md", "", "", "e75892d4-044b-4d09-8abf-855ed5fa4f20.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cd76318f29306b8e90dd81f65ec88bd3a309cc5e/e2e/7c27a539-13f8-415c-9532-2e5a2ec12fb6.md", "", "", "7c27a539-13f8-415c-9532-2e5a2ec12fb6.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cd76318f29306b8e90dd81f65ec88bd3a309cc5e/e2e/acfcdcf1-9169-4293-80a2-ed8587dc1452.md", "", "", "acfcdcf1-9169-4293-80a2-ed8587dc1452.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d92d0733590fd5894259fe924c051e691cdc04b6/e2e/3414b3fb-bece-4317-851b-2080f8fb27f6.md", "", "", "3414b3fb-bece-4317-851b-2080f8fb27f6.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d92d0733590fd5894259fe924c051e691cdc04b6/e2e/7f880207-dc89-446c-99de-d3b10383c0e3.md", "", "", "7f880207-dc89-446c-99de-d3b10383c0e3.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d92d0733590fd5894259fe924c051e691cdc04b6/e2e/f5fea1db-e187-4e4e-a926-19f06f96b483.md", "", "", "f5fea1db-e187-4e4e-a926-19f06f96b483.md") | Out-Null

Style-Hyperlink $wsZhCn.Range("A2")
Style-Hyperlink $wsZhCn.Range("I2")
Style-Hyperlink $wsZhCn.Range("A3")
Style-Hyperlink $wsZhCn.Range("A4")
Style-Hyperlink $wsZhCn.Range("A5")
Style-Hyperlink $wsZhCn.Range("A6")
Style-Hyperlink $wsZhCn.Range("A7")

# Re-apply the date format to the other already-existing date cells that
# Hyperlinks.Delete()/Add() may have reset the font on (H/K columns untouched,
# but make sure rows 2-4 keep their number format).
Style-DateText $wsZhCn.Range("H2")
Style-DateText $wsZhCn.Range("K2")
Style-DateText $wsZhCn.Range("H3")
Style-DateText $wsZhCn.Range("K3")
Style-DateText $wsZhCn.Range("H4")
Style-DateText $wsZhCn.Range("K4")

$tblZhCn = $wsZhCn.ListObjects.Item(1)
$tblZhCn.Resize($wsZhCn.Range("A1:P7"))

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

# Row 5: new 3414b3fb entry
$wsDeDe.Range("A5").Value = "3414b3fb-bece-4317-851b-2080f8fb27f6.md"
$wsDeDe.Range("B5").Value = ".md"
$wsDeDe.Range("C5").Value = "Ready for handoff"
$wsDeDe.Range("D5").Value = "e2e"
$wsDeDe.Range("E5").Value = "ht"
$wsDeDe.Range("F5").Value = "False"
$wsDeDe.Range("G5").Value = "3414b3fb-bece-4317-851b-2080f8fb27f6.9a7d471b31670a0c41801e2ea85158d6511dc541.de-de.xlf"
$wsDeDe.Range("H5").Value = "2016-08-31 12:51:36"
$wsDeDe.Range("K5").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("M5").Value = "True"
$wsDeDe.Range("O5").Value = "False"
Style-DateText $wsDeDe.Range("H5")
Style-DateText $wsDeDe.Range("K5")

# Row 6: previous 7f880207 entry, shifted down
$wsDeDe.Range("A6").Value = "7f880207-dc89-446c-99de-d3b10383c0e3.md"
$wsDeDe.Range("B6").Value = ".md"
$wsDeDe.Range("C6").Value = "Ready for handoff"
$wsDeDe.Range("D6").Value = "e2e"
$wsDeDe.Range("E6").Value = "ht"
$wsDeDe.Range("F6").Value = "False"
$wsDeDe.Range("G6").Value = "7f880207-dc89-446c-99de-d3b10383c0e3.fc6f83bef6d89ff41f4a184c043c9e87ea8c1cdf.de-de.xlf"
$wsDeDe.Range("H6").Value = "2016-08-31 12:49:20"
$wsDeDe.Range("K6").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("M6").Value = "True"
$wsDeDe.Range("O6").Value = "False"
Style-DateText $wsDeDe.Range("H6")
Style-DateText $wsDeDe.Range("K6")

# Row 7: new f5fea1db entry
$wsDeDe.Range("A7").Value = "f5fea1db-e187-4e4e-a926-19f06f96b483.md"
$wsDeDe.Range("B7").Value = ".md"
$wsDeDe.Range("C7").Value = "Ready for handoff"
$wsDeDe.Range("D7").Value = "e2e"
$wsDeDe.Range("E7").Value = "ht"
$wsDeDe.Range("F7").Value = "False"
$wsDeDe.Range("G7").Value = "f5fea1db-e187-4e4e-a926-19f06f96b483.9e74bf257fb01b4211733f5bb469e0eb170f6bfd.de-de.xlf"
$wsDeDe.Range("H7").Value = "2016-08-31 12:51:36"
$wsDeDe.Range("K7").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("M7").Value = "True"
$wsDeDe.Range("O7").Value = "False"
Style-DateText $wsDeDe.Range("H7")
Style-DateText $wsDeDe.Range("K7")

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/90e4f0ebcb4243bb0e92ca46307a9286e675fe81/e2e/e75892d4-044b-4d09-8abf-855ed5fa4f20.md", "", "", "e75892d4-044b-4d09-8abf-855ed5fa4f20.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/79a06fa96bbeca38e875708771a99b54a4486f16/e2e/e75892d4-044b-4d09-8abf-855ed5fa4f20.md", "", "", "e75892d4-044b-4d09-8abf-855ed5fa4f20.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cd76318f29306b8e90dd81f65ec88bd3a309cc5e/e2e/7c27a539-13f8-415c-9532-2e5a2ec12fb6.md", "", "", "7c27a539-13f8-415c-9532-2e5a2ec12fb6.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cd76318f29306b8e90dd81f65ec88bd3a309cc5e/e2e/acfcdcf1-9169-4293-80a2-ed8587dc1452.md", "", "", "acfcdcf1-9169-4293-80a2-ed8587dc1452.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d92d0733590fd5894259fe924c051e691cdc04b6/e2e/3414b3fb-bece-4317-851b-2080f8fb27f6.md", "", "", "3414b3fb-bece-4317-851b-2080f8fb27f6.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d92d0733590fd5894259fe924c051e691cdc04b6/e2e/7f880207-dc89-446c-99de-d3b10383c0e3.md", "", "", "7f880207-dc89-446c-99de-d3b10383c0e3.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d92d0733590fd5894259fe924c051e691cdc04b6/e2e/f5fea1db-e187-4e4e-a926-19f06f96b483.md", "", "", "f5fea1db-e187-4e4e-a926-19f06f96b483.md") | Out-Null

Style-Hyperlink $wsDeDe.Range("A2")
Style-Hyperlink $wsDeDe.Range("I2")
Style-Hyperlink $wsDeDe.Range("A3")
Style-Hyperlink $wsDeDe.Range("A4")
Style-Hyperlink $wsDeDe.Range("A5")
Style-Hyperlink $wsDeDe.Range("A6")
Style-Hyperlink $wsDeDe.Range("A7")

Style-DateText $wsDeDe.Range("H2")
Style-DateText $wsDeDe.Range("K2")
Style-DateText $wsDeDe.Range("H3")
Style-DateText $wsDeDe.Range("K3")
Style-DateText $wsDeDe.Range("H4")
Style-DateText $wsDeDe.Range("K4")

$tblDeDe = $wsDeDe.ListObjects.Item(1)
$tblDeDe.Resize($wsDeDe.Range("A1:P7"))
